$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block: card holder name / account number ---
$ws.Range("C2").Value = "Hartmut"

# B3 holds the card number as TEXT (not a number) in the original file.
# A plain Value assignment of an all-digit string auto-converts to a
# number, so instead stage the literal text in a scratch cell far outside
# the used range, copy it across as a value (preserving its text type and
# B3's existing cell style), then wipe the scratch cell.
$ws.Range("ZZ1").Formula = "=""2570314725427075"""
$ws.Range("ZZ1").Copy()
$ws.Range("B3").PasteSpecial(-4163)
$ws.Range("ZZ1").Value = ""

$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance line ---
$ws.Range("D5").Value = "KONTOSTAND AM 03.07.2024"

# --- Row 6 ---
$ws.Range("B6").Value = "07.07."
$ws.Range("C6").Value = "08.07."
$ws.Range("D6").Value = "KARTENZAHLUNG ARAL TANKSTELLE"
$ws.Range("E6").Value = "69,23-"

# --- Row 7 ---
$ws.Range("B7").Value = "11.07."
$ws.Range("C7").Value = "12.07."
$ws.Range("D7").Value = "KARTENZAHLUNG JET TANKSTELLE"
$ws.Range("E7").Value = "68,00-"

# --- Row 8 ---
$ws.Range("B8").Value = "15.07."
$ws.Range("C8").Value = "16.07."
$ws.Range("D8").Value = "BEITRAG Allianz SE K-68381771"
$ws.Range("E8").Value = "54,65-"

# --- Row 9 ---
$ws.Range("B9").Value = "16.07."
$ws.Range("C9").Value = "17.07."
$ws.Range("D9").Value = "MCDONALDS Mellrichstadt"
$ws.Range("E9").Value = "14,53-"

# --- Row 10: the former 5th transaction is removed, row becomes blank
#     (matching the blank template row 11 in both content and formatting) ---
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = ""
$ws.Range("E11").Copy()
$ws.Range("E10").PasteSpecial(-4122)

# --- Closing balance line ---
$ws.Range("D12").Value = "KONTOSTAND AM 21.07.2024"
$ws.Range("E12").Value = "206,41-"

# --- Next statement date ---
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 30.07.2024"
